$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 734.86365
$ws.Range("I15").Value = 734.86365
$ws.Range("K15").Value = 2204.59095
$ws.Range("M15").Value = -2035.59095
$ws.Range("H40").Value = 10769.533
$ws.Range("I40").Value = 16529.143
$ws.Range("J40").Value = 5729.875
$ws.Range("K40").Value = 16529.143
$ws.Range("L40").Value = 5729.875
$ws.Range("M40").Value = -16354.143
$ws.Range("N40").Value = -6079.875
$ws.Range("H69").Value = 28579184
$ws.Range("I69").Value = 7980
$ws.Range("J69").Value = 29419514
$ws.Range("K69").Value = 23940
$ws.Range("L69").Value = 88258542
$ws.Range("M69").Value = -23066
$ws.Range("N69").Value = -88260290
$ws.Range("H72").Value = 28579184
$ws.Range("I72").Value = 7980
$ws.Range("J72").Value = 29419514
$ws.Range("K72").Value = 71820
$ws.Range("L72").Value = 264775626
$ws.Range("M72").Value = -67452
$ws.Range("N72").Value = -264784362
$ws.Range("H74").Value = 6170.875
$ws.Range("I74").Value = 5563.6665
$ws.Range("J74").Value = 7992.5
$ws.Range("K74").Value = 5563.6665
$ws.Range("L74").Value = 7992.5
$ws.Range("M74").Value = -4627.6665
$ws.Range("N74").Value = -9864.5
$ws.Range("H77").Value = 6170.875
$ws.Range("I77").Value = 5563.6665
$ws.Range("J77").Value = 7992.5
$ws.Range("K77").Value = 27818.3325
$ws.Range("L77").Value = 39962.5
$ws.Range("M77").Value = -23138.3325
$ws.Range("N77").Value = -49322.5
$ws.Range("H82").Value = 1557.5555
$ws.Range("I82").Value = 1557.5555
$ws.Range("K82").Value = 4672.666499999999
$ws.Range("M82").Value = -4266.666499999999
$ws.Range("H85").Value = 1557.5555
$ws.Range("I85").Value = 1557.5555
$ws.Range("K85").Value = 4672.666499999999
$ws.Range("M85").Value = -3268.666499999999
$ws.Range("H112").Value = 11840.667
$ws.Range("I112").Value = 2985
$ws.Range("J112").Value = 12645.728
$ws.Range("K112").Value = 8955
$ws.Range("L112").Value = 37937.18399999999
$ws.Range("M112").Value = -7847
$ws.Range("N112").Value = -40153.18399999999
$ws.Range("H113").Value = 7537.9
$ws.Range("I113").Value = 7302
$ws.Range("J113").Value = 7596.875
$ws.Range("K113").Value = 7302
$ws.Range("L113").Value = 7596.875
$ws.Range("M113").Value = -4048
$ws.Range("N113").Value = -14104.875
$ws.Range("H118").Value = 10417089
$ws.Range("I118").Value = 11905112
$ws.Range("J118").Value = 928
$ws.Range("K118").Value = 35715336
$ws.Range("L118").Value = 2784
$ws.Range("M118").Value = -35713679
$ws.Range("N118").Value = -6098
$ws.Range("H137").Value = 138813.08
$ws.Range("I137").Value = 357834.8
$ws.Range("J137").Value = 1924.5
$ws.Range("K137").Value = 1073504.4
$ws.Range("L137").Value = 5773.5
$ws.Range("M137").Value = -1070954.4
$ws.Range("N137").Value = -10873.5
$ws.Range("H138").Value = 7193.706
$ws.Range("I138").Value = 10000
$ws.Range("J138").Value = 7018.3125
$ws.Range("K138").Value = 30000
$ws.Range("L138").Value = 21054.9375
$ws.Range("M138").Value = -24860
$ws.Range("N138").Value = -31334.9375

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8099897
$ws.Range("I45").Value = 9617253
$ws.Range("K45").Value = 9617253
$ws.Range("M45").Value = -9616876
$ws.Range("H110").Value = 3087890.5
$ws.Range("I110").Value = 3473727
$ws.Range("J110").Value = 1200
$ws.Range("K110").Value = 3473727
$ws.Range("L110").Value = 1200
$ws.Range("M110").Value = -3471682
$ws.Range("N110").Value = -5290
$ws.Range("H122").Value = 1100270.4
$ws.Range("I122").Value = 3598
$ws.Range("J122").Value = 1897850.1
$ws.Range("K122").Value = 10794
$ws.Range("L122").Value = 5693550.300000001
$ws.Range("M122").Value = -8344
$ws.Range("N122").Value = -5698450.300000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3228.89
$ws.Range("I31").Value = 1740.3529
$ws.Range("J31").Value = 3533.771
$ws.Range("K31").Value = 1740.3529
$ws.Range("L31").Value = 3533.771
$ws.Range("M31").Value = -1445.3529
$ws.Range("N31").Value = -4123.771000000001
$ws.Range("H34").Value = 3228.89
$ws.Range("I34").Value = 1740.3529
$ws.Range("J34").Value = 3533.771
$ws.Range("K34").Value = 1740.3529
$ws.Range("L34").Value = 3533.771
$ws.Range("M34").Value = -1538.3529
$ws.Range("N34").Value = -3937.771
$ws.Range("H107").Value = 34483524
$ws.Range("I107").Value = 800.8333
$ws.Range("J107").Value = 200000600
$ws.Range("K107").Value = 800.8333
$ws.Range("L107").Value = 200000600
$ws.Range("M107").Value = 1119.1667
$ws.Range("N107").Value = -200004440

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3137.9722
$ws.Range("I68").Value = 2893.5789
$ws.Range("J68").Value = 3411.1177
$ws.Range("K68").Value = 8680.736699999999
$ws.Range("L68").Value = 10233.3531
$ws.Range("M68").Value = -7869.736699999999
$ws.Range("N68").Value = -11855.3531
$ws.Range("H71").Value = 3137.9722
$ws.Range("I71").Value = 2893.5789
$ws.Range("J71").Value = 3411.1177
$ws.Range("K71").Value = 26042.2101
$ws.Range("L71").Value = 30700.0593
$ws.Range("M71").Value = -21986.2101
$ws.Range("N71").Value = -38812.05929999999
$ws.Range("H92").Value = 788.2
$ws.Range("I92").Value = 712.75
$ws.Range("J92").Value = 874.4286
$ws.Range("K92").Value = 2138.25
$ws.Range("L92").Value = 2623.2858
$ws.Range("M92").Value = -890.25
$ws.Range("N92").Value = -5119.2858
$ws.Range("H113").Value = 2516.0303
$ws.Range("I113").Value = 4055.4443
$ws.Range("J113").Value = 1938.75
$ws.Range("K113").Value = 12166.3329
$ws.Range("L113").Value = 5816.25
$ws.Range("M113").Value = -9996.332900000001
$ws.Range("N113").Value = -10156.25
$ws.Range("H131").Value = 2909.3333
$ws.Range("J131").Value = 4249.75
$ws.Range("L131").Value = 12749.25
$ws.Range("N131").Value = -22829.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 621.6875
$ws.Range("I107").Value = 703.7692
$ws.Range("J107").Value = 266
$ws.Range("K107").Value = 703.7692
$ws.Range("L107").Value = 266
$ws.Range("M107").Value = 1216.2308
$ws.Range("N107").Value = -4106
$ws.Range("H113").Value = 6947373
$ws.Range("I113").Value = 11906580
$ws.Range("J113").Value = 4483.3
$ws.Range("K113").Value = 11906580
$ws.Range("L113").Value = 4483.3
$ws.Range("M113").Value = -11904410
$ws.Range("N113").Value = -8823.299999999999
$ws.Range("H122").Value = 809253.25
$ws.Range("I122").Value = 890078.6
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2670235.8
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -2667785.8
$ws.Range("N122").Value = -7900

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 26985.756
$ws.Range("I136").Value = 34307.418
$ws.Range("J136").Value = 4288.6
$ws.Range("K136").Value = 102922.254
$ws.Range("L136").Value = 12865.8
$ws.Range("M136").Value = -100372.254
$ws.Range("N136").Value = -17965.8
$ws.Range("H141").Value = 118271.664
$ws.Range("J141").Value = 118271.664
$ws.Range("L141").Value = 118271.664
$ws.Range("N141").Value = -128631.664

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5843.3184
$ws.Range("I136").Value = 7389.516
$ws.Range("J136").Value = 2156.2307
$ws.Range("K136").Value = 22168.548
$ws.Range("L136").Value = 6468.6921
$ws.Range("M136").Value = -19618.548
$ws.Range("N136").Value = -11568.6921
